$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 914.2857
$ws.Range("I18").Value = 914.2857
$ws.Range("K18").Value = 914.2857
$ws.Range("M18").Value = -630.2857

$ws.Range("H33").Value = 235.6923
$ws.Range("I33").Value = 243.09091
$ws.Range("K33").Value = 243.09091
$ws.Range("M33").Value = -14.09091000000001

$ws.Range("H43").Value = 5280.5713
$ws.Range("J43").Value = 6772.3887
$ws.Range("L43").Value = 6772.3887
$ws.Range("N43").Value = -6910.3887

$ws.Range("H61").Value = 1244.875
$ws.Range("I61").Value = 191.8
$ws.Range("K61").Value = 575.4000000000001
$ws.Range("M61").Value = -403.4000000000001

$ws.Range("H80").Value = 980.5789
$ws.Range("J80").Value = 1798.5
$ws.Range("L80").Value = 5395.5
$ws.Range("N80").Value = -7391.5

$ws.Range("H83").Value = 980.5789
$ws.Range("J83").Value = 1798.5
$ws.Range("L83").Value = 16186.5
$ws.Range("N83").Value = -26170.5

$ws.Range("H86").Value = 5499.6665
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 5499.6665
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 5499.6665
$ws.Range("N86").Value = -7745.6665

$ws.Range("H89").Value = 5499.6665
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 5499.6665
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 27498.3325
$ws.Range("N89").Value = -38730.3325

$ws.Range("H100").Value = 1538.2258
$ws.Range("J100").Value = 2832.4546
$ws.Range("L100").Value = 2832.4546
$ws.Range("N100").Value = -3914.4546

$ws.Range("H112").Value = 1896.4706
$ws.Range("J112").Value = 2385.6365
$ws.Range("L112").Value = 7156.9095
$ws.Range("N112").Value = -9372.9095

$ws.Range("H136").Value = 97000
$ws.Range("J136").Value = 97000
$ws.Range("L136").Value = 97000
$ws.Range("N136").Value = -107200

$ws.Range("H137").Value = 3333.53
$ws.Range("I137").Value = 2030.0984
$ws.Range("K137").Value = 6090.2952
$ws.Range("M137").Value = -3540.2952

$ws.Range("H138").Value = 4991.566
$ws.Range("J138").Value = 5003.3267
$ws.Range("L138").Value = 15009.9801
$ws.Range("N138").Value = -25289.9801

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1770569
$ws.Range("I2").Value = 1958816.8
$ws.Range("K2").Value = 1958816.8
$ws.Range("M2").Value = -1958703.8

$ws.Range("H32").Value = 5109.62
$ws.Range("I32").Value = 2907.186
$ws.Range("K32").Value = 2907.186
$ws.Range("M32").Value = -2620.186

$ws.Range("H33").Value = 11604.833
$ws.Range("I33").Value = 8720
$ws.Range("J33").Value = 26029
$ws.Range("K33").Value = 8720
$ws.Range("L33").Value = 26029
$ws.Range("M33").Value = -8391
$ws.Range("N33").Value = -26687

$ws.Range("H36").Value = 1431725.9
$ws.Range("I36").Value = 2610.4
$ws.Range("K36").Value = 2610.4
$ws.Range("M36").Value = -2264.4

$ws.Range("H116").Value = 1770569
$ws.Range("I116").Value = 1958816.8
$ws.Range("K116").Value = 1958816.8
$ws.Range("M116").Value = -1956522.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1770569
$ws.Range("I3").Value = 1958816.8
$ws.Range("K3").Value = 1958816.8
$ws.Range("M3").Value = -1958702.8

$ws.Range("H105").Value = 111114000
$ws.Range("I105").Value = 166668670
$ws.Range("J105").Value = 4673.6665
$ws.Range("K105").Value = 166668670
$ws.Range("L105").Value = 4673.6665
$ws.Range("M105").Value = -166666923
$ws.Range("N105").Value = -8167.6665

$ws.Range("H134").Value = 4262.0186
$ws.Range("I134").Value = 3026.1135
$ws.Range("K134").Value = 9078.3405
$ws.Range("M134").Value = -6543.3405

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31254806
$ws.Range("I31").Value = 50002388
$ws.Range("K31").Value = 50002388
$ws.Range("M31").Value = -50002093

$ws.Range("H34").Value = 31254806
$ws.Range("I34").Value = 50002388
$ws.Range("K34").Value = 50002388
$ws.Range("M34").Value = -50002186

$ws.Range("H58").Value = 4593.8945
$ws.Range("I58").Value = 2483.5454
$ws.Range("J58").Value = 7495.625
$ws.Range("K58").Value = 2483.5454
$ws.Range("L58").Value = 7495.625
$ws.Range("M58").Value = -2280.5454
$ws.Range("N58").Value = -7901.625

$ws.Range("H60").Value = 25117
$ws.Range("J60").Value = 29814.715
$ws.Range("L60").Value = 29814.715
$ws.Range("N60").Value = -30836.715

$ws.Range("H94").Value = 1383.4706
$ws.Range("I94").Value = 1665.6
$ws.Range("K94").Value = 1665.6
$ws.Range("M94").Value = -1214.6

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws.Range("H136").Value = 4593.8945
$ws.Range("I136").Value = 2483.5454
$ws.Range("J136").Value = 7495.625
$ws.Range("K136").Value = 7450.6362
$ws.Range("L136").Value = 22486.875
$ws.Range("M136").Value = -4900.6362
$ws.Range("N136").Value = -27586.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1399
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1399
$ws.Range("K20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").Value = 4197
$ws.Range("N20").Value = -4651

$ws.Range("H25").Value = 527.7143
$ws.Range("J25").Value = 349.5
$ws.Range("L25").Value = 1048.5
$ws.Range("N25").Value = -1386.5

$ws.Range("H30").Value = 527.7143
$ws.Range("J30").Value = 349.5
$ws.Range("L30").Value = 1048.5
$ws.Range("N30").Value = -1252.5

$ws.Range("H118").Value = 1676.3334
$ws.Range("I118").Value = 1676.3334
$ws.Range("K118").Value = 5029.0002
$ws.Range("M118").Value = -3786.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11999.6
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 10000
$ws.Range("M70").Value = -9730

$ws.Range("H73").Value = 11999.6
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 10000
$ws.Range("M73").Value = -9064

$ws.Range("H80").Value = 5866.65
$ws.Range("I80").Value = 7419.7144
$ws.Range("J80").Value = 5030.385
$ws.Range("K80").Value = 7419.7144
$ws.Range("L80").Value = 5030.385
$ws.Range("M80").Value = -6421.7144
$ws.Range("N80").Value = -7026.385

$ws.Range("H83").Value = 5866.65
$ws.Range("I83").Value = 7419.7144
$ws.Range("J83").Value = 5030.385
$ws.Range("K83").Value = 37098.572
$ws.Range("L83").Value = 25151.925
$ws.Range("M83").Value = -32106.572
$ws.Range("N83").Value = -35135.925

$ws.Range("H97").Value = 43479344
$ws.Range("I97").Value = 837.26666
$ws.Range("J97").Value = 125001544
$ws.Range("K97").Value = 837.26666
$ws.Range("L97").Value = 125001544
$ws.Range("M97").Value = -341.26666
$ws.Range("N97").Value = -125002536

$ws.Range("H102").Value = 2155.6667
$ws.Range("I102").Value = 2057.2856
$ws.Range("K102").Value = 2057.2856
$ws.Range("M102").Value = -435.2856000000002

$ws.Range("H113").Value = 3457.0715
$ws.Range("I113").Value = 2765.75
$ws.Range("J113").Value = 4378.8335
$ws.Range("K113").Value = 2765.75
$ws.Range("L113").Value = 4378.8335
$ws.Range("M113").Value = -595.75
$ws.Range("N113").Value = -8718.833500000001

$ws.Range("H126").Value = 3414.125
$ws.Range("I126").Value = 2340.25
$ws.Range("K126").Value = 7020.75
$ws.Range("M126").Value = -4550.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 986.2222
$ws.Range("I16").Value = 986.2222
$ws.Range("K16").Value = 986.2222
$ws.Range("M16").Value = -816.2222

$ws.Range("H22").Value = 4434.0835
$ws.Range("I22").Value = 2785.4285
$ws.Range("K22").Value = 2785.4285
$ws.Range("M22").Value = -2490.4285

$ws.Range("H27").Value = 4434.0835
$ws.Range("I27").Value = 2785.4285
$ws.Range("K27").Value = 2785.4285
$ws.Range("M27").Value = -2678.4285

$ws.Range("H61").Value = 2433.5625
$ws.Range("I61").Value = 2670.8928
$ws.Range("J61").Value = 772.25
$ws.Range("K61").Value = 2670.8928
$ws.Range("L61").Value = 772.25
$ws.Range("M61").Value = -2468.8928
$ws.Range("N61").Value = -1176.25

$ws.Range("H93").Value = 3372.875
$ws.Range("I93").Value = 2531.5
$ws.Range("K93").Value = 2531.5
$ws.Range("M93").Value = -1283.5

$ws.Range("H113").Value = 2433.5625
$ws.Range("I113").Value = 2670.8928
$ws.Range("J113").Value = 772.25
$ws.Range("K113").Value = 2670.8928
$ws.Range("L113").Value = 772.25
$ws.Range("M113").Value = -500.8928000000001
$ws.Range("N113").Value = -5112.25

$ws.Range("H132").Value = 5355
$ws.Range("I132").Value = 4678.2905
$ws.Range("J132").Value = 6403.9
$ws.Range("K132").Value = 14034.8715
$ws.Range("L132").Value = 19211.7
$ws.Range("M132").Value = -11504.8715
$ws.Range("N132").Value = -24271.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3134.5386
$ws.Range("I122").Value = 2818.0908
$ws.Range("K122").Value = 8454.2724
$ws.Range("M122").Value = -6004.2724

$ws.Range("H132").Value = 6267.45
$ws.Range("I132").Value = 4811.1113
$ws.Range("J132").Value = 7459
$ws.Range("K132").Value = 14433.3339
$ws.Range("L132").Value = 22377
$ws.Range("M132").Value = -11903.3339
$ws.Range("N132").Value = -27437
